$d = $word.ActiveDocument

$d.Content.Find.Execute("60+13=", $true, $true, $false, $false, $false, $true, 1, $false, "5+31=", 2) | Out-Null
$d.Content.Find.Execute("7+65=", $true, $true, $false, $false, $false, $true, 1, $false, "72+7=", 2) | Out-Null
$d.Content.Find.Execute("16+77=", $true, $true, $false, $false, $false, $true, 1, $false, "71+13=", 2) | Out-Null
$d.Content.Find.Execute("96-90=", $true, $true, $false, $false, $false, $true, 1, $false, "68-24=", 2) | Out-Null
$d.Content.Find.Execute("35+29=", $true, $true, $false, $false, $false, $true, 1, $false, "39+12=", 2) | Out-Null
$d.Content.Find.Execute("2+16=", $true, $true, $false, $false, $false, $true, 1, $false, "81-18=", 2) | Out-Null
$d.Content.Find.Execute("27+15=", $true, $true, $false, $false, $false, $true, 1, $false, "51-4=", 2) | Out-Null
$d.Content.Find.Execute("45-17=", $true, $true, $false, $false, $false, $true, 1, $false, "0+34=", 2) | Out-Null
$d.Content.Find.Execute("44-38=", $true, $true, $false, $false, $false, $true, 1, $false, "8+48=", 2) | Out-Null
$d.Content.Find.Execute("61-56=", $true, $true, $false, $false, $false, $true, 1, $false, "23+72=", 2) | Out-Null
$d.Content.Find.Execute("5+39=", $true, $true, $false, $false, $false, $true, 1, $false, "4+92=", 2) | Out-Null
$d.Content.Find.Execute("22+57=", $true, $true, $false, $false, $false, $true, 1, $false, "5+83=", 2) | Out-Null
$d.Content.Find.Execute("37+55=", $true, $true, $false, $false, $false, $true, 1, $false, "4+14=", 2) | Out-Null
$d.Content.Find.Execute("65+0=", $true, $true, $false, $false, $false, $true, 1, $false, "26-14=", 2) | Out-Null
$d.Content.Find.Execute("75-72=", $true, $true, $false, $false, $false, $true, 1, $false, "77-44=", 2) | Out-Null
$d.Content.Find.Execute("10+51=", $true, $true, $false, $false, $false, $true, 1, $false, "13+0=", 2) | Out-Null
$d.Content.Find.Execute("14+38=", $true, $true, $false, $false, $false, $true, 1, $false, "28+38=", 2) | Out-Null
$d.Content.Find.Execute("41+6=", $true, $true, $false, $false, $false, $true, 1, $false, "16+48=", 2) | Out-Null
$d.Content.Find.Execute("82+7=", $true, $true, $false, $false, $false, $true, 1, $false, "18-15=", 2) | Out-Null
$d.Content.Find.Execute("79-0=", $true, $true, $false, $false, $false, $true, 1, $false, "26+62=", 2) | Out-Null
$d.Content.Find.Execute("45+12=", $true, $true, $false, $false, $false, $true, 1, $false, "66-29=", 2) | Out-Null
$d.Content.Find.Execute("54+35=", $true, $true, $false, $false, $false, $true, 1, $false, "1+23=", 2) | Out-Null
$d.Content.Find.Execute("32+2=", $true, $true, $false, $false, $false, $true, 1, $false, "21+63=", 2) | Out-Null
$d.Content.Find.Execute("77-53=", $true, $true, $false, $false, $false, $true, 1, $false, "28-21=", 2) | Out-Null
$d.Content.Find.Execute("58+27=", $true, $true, $false, $false, $false, $true, 1, $false, "88-11=", 2) | Out-Null
$d.Content.Find.Execute("21+4=", $true, $true, $false, $false, $false, $true, 1, $false, "66+26=", 2) | Out-Null
$d.Content.Find.Execute("49+10=", $true, $true, $false, $false, $false, $true, 1, $false, "37-5=", 2) | Out-Null
$d.Content.Find.Execute("17+16=", $true, $true, $false, $false, $false, $true, 1, $false, "54+11=", 2) | Out-Null
$d.Content.Find.Execute("67-58=", $true, $true, $false, $false, $false, $true, 1, $false, "16+61=", 2) | Out-Null
$d.Content.Find.Execute("46+30=", $true, $true, $false, $false, $false, $true, 1, $false, "91-12=", 2) | Out-Null
$d.Content.Find.Execute("70-33=", $true, $true, $false, $false, $false, $true, 1, $false, "60-40=", 2) | Out-Null
$d.Content.Find.Execute("3+76=", $true, $true, $false, $false, $false, $true, 1, $false, "59+12=", 2) | Out-Null
$d.Content.Find.Execute("64+32=", $true, $true, $false, $false, $false, $true, 1, $false, "53-43=", 2) | Out-Null
$d.Content.Find.Execute("72-20=", $true, $true, $false, $false, $false, $true, 1, $false, "81-9=", 2) | Out-Null
$d.Content.Find.Execute("92-50=", $true, $true, $false, $false, $false, $true, 1, $false, "91-73=", 2) | Out-Null
$d.Content.Find.Execute("2+79=", $true, $true, $false, $false, $false, $true, 1, $false, "29+29=", 2) | Out-Null
$d.Content.Find.Execute("5+2=", $true, $true, $false, $false, $false, $true, 1, $false, "55+30=", 2) | Out-Null
$d.Content.Find.Execute("22-6=", $true, $true, $false, $false, $false, $true, 1, $false, "35+10=", 2) | Out-Null
$d.Content.Find.Execute("47-22=", $true, $true, $false, $false, $false, $true, 1, $false, "80-77=", 2) | Out-Null
$d.Content.Find.Execute("88-86=", $true, $true, $false, $false, $false, $true, 1, $false, "91-78=", 2) | Out-Null
$d.Content.Find.Execute("22+68=", $true, $true, $false, $false, $false, $true, 1, $false, "29+26=", 2) | Out-Null
$d.Content.Find.Execute("13+26=", $true, $true, $false, $false, $false, $true, 1, $false, "34+11=", 2) | Out-Null
$d.Content.Find.Execute("95-60=", $true, $true, $false, $false, $false, $true, 1, $false, "2+96=", 2) | Out-Null
$d.Content.Find.Execute("7+70=", $true, $true, $false, $false, $false, $true, 1, $false, "1+14=", 2) | Out-Null
$d.Content.Find.Execute("93-88=", $true, $true, $false, $false, $false, $true, 1, $false, "87+6=", 2) | Out-Null
$d.Content.Find.Execute("4+77=", $true, $true, $false, $false, $false, $true, 1, $false, "8-5=", 2) | Out-Null
$d.Content.Find.Execute("94-36=", $true, $true, $false, $false, $false, $true, 1, $false, "85+13=", 2) | Out-Null
$d.Content.Find.Execute("52+38=", $true, $true, $false, $false, $false, $true, 1, $false, "93-7=", 2) | Out-Null
$d.Content.Find.Execute("75-29=", $true, $true, $false, $false, $false, $true, 1, $false, "52-3=", 2) | Out-Null
$d.Content.Find.Execute("93-28=", $true, $true, $false, $false, $false, $true, 1, $false, "21+70=", 2) | Out-Null
$d.Content.Find.Execute("29+5=", $true, $true, $false, $false, $false, $true, 1, $false, "7+50=", 2) | Out-Null
$d.Content.Find.Execute("18+70=", $true, $true, $false, $false, $false, $true, 1, $false, "16+23=", 2) | Out-Null
$d.Content.Find.Execute("52-48=", $true, $true, $false, $false, $false, $true, 1, $false, "47+20=", 2) | Out-Null
$d.Content.Find.Execute("12+19=", $true, $true, $false, $false, $false, $true, 1, $false, "31-15=", 2) | Out-Null
$d.Content.Find.Execute("9+35=", $true, $true, $false, $false, $false, $true, 1, $false, "12+34=", 2) | Out-Null
$d.Content.Find.Execute("44-3=", $true, $true, $false, $false, $false, $true, 1, $false, "34-6=", 2) | Out-Null
$d.Content.Find.Execute("95-7=", $true, $true, $false, $false, $false, $true, 1, $false, "89-26=", 2) | Out-Null
$d.Content.Find.Execute("65-61=", $true, $true, $false, $false, $false, $true, 1, $false, "96-3=", 2) | Out-Null
$d.Content.Find.Execute("85+2=", $true, $true, $false, $false, $false, $true, 1, $false, "44-34=", 2) | Out-Null
$d.Content.Find.Execute("98-98=", $true, $true, $false, $false, $false, $true, 1, $false, "30+7=", 2) | Out-Null
$d.Content.Find.Execute("89-41=", $true, $true, $false, $false, $false, $true, 1, $false, "80-76=", 2) | Out-Null
$d.Content.Find.Execute("30+2=", $true, $true, $false, $false, $false, $true, 1, $false, "61-58=", 2) | Out-Null
$d.Content.Find.Execute("19+37=", $true, $true, $false, $false, $false, $true, 1, $false, "52-31=", 2) | Out-Null
$d.Content.Find.Execute("99-44=", $true, $true, $false, $false, $false, $true, 1, $false, "16+71=", 2) | Out-Null
$d.Content.Find.Execute("4+94=", $true, $true, $false, $false, $false, $true, 1, $false, "10+68=", 2) | Out-Null
$d.Content.Find.Execute("3+58=", $true, $true, $false, $false, $false, $true, 1, $false, "32+21=", 2) | Out-Null
$d.Content.Find.Execute("10+55=", $true, $true, $false, $false, $false, $true, 1, $false, "17+17=", 2) | Out-Null
$d.Content.Find.Execute("77-64=", $true, $true, $false, $false, $false, $true, 1, $false, "76-42=", 2) | Out-Null
$d.Content.Find.Execute("94-20=", $true, $true, $false, $false, $false, $true, 1, $false, "64-12=", 2) | Out-Null
$d.Content.Find.Execute("45+30=", $true, $true, $false, $false, $false, $true, 1, $false, "17+10=", 2) | Out-Null
$d.Content.Find.Execute("16+7=", $true, $true, $false, $false, $false, $true, 1, $false, "94-19=", 2) | Out-Null
$d.Content.Find.Execute("99-38=", $true, $true, $false, $false, $false, $true, 1, $false, "48-23=", 2) | Out-Null
$d.Content.Find.Execute("96-41=", $true, $true, $false, $false, $false, $true, 1, $false, "23+9=", 2) | Out-Null
$d.Content.Find.Execute("40+13=", $true, $true, $false, $false, $false, $true, 1, $false, "61-34=", 2) | Out-Null
$d.Content.Find.Execute("55+13=", $true, $true, $false, $false, $false, $true, 1, $false, "39-3=", 2) | Out-Null
$d.Content.Find.Execute("32+51=", $true, $true, $false, $false, $false, $true, 1, $false, "74-48=", 2) | Out-Null
$d.Content.Find.Execute("73-59=", $true, $true, $false, $false, $false, $true, 1, $false, "87+11=", 2) | Out-Null
$d.Content.Find.Execute("88-30=", $true, $true, $false, $false, $false, $true, 1, $false, "94-12=", 2) | Out-Null
$d.Content.Find.Execute("26+13=", $true, $true, $false, $false, $false, $true, 1, $false, "29+59=", 2) | Out-Null
$d.Content.Find.Execute("51-18=", $true, $true, $false, $false, $false, $true, 1, $false, "7+43=", 2) | Out-Null
$d.Content.Find.Execute("36+40=", $true, $true, $false, $false, $false, $true, 1, $false, "36+34=", 2) | Out-Null
$d.Content.Find.Execute("11+1=", $true, $true, $false, $false, $false, $true, 1, $false, "35+49=", 2) | Out-Null
$d.Content.Find.Execute("36+62=", $true, $true, $false, $false, $false, $true, 1, $false, "23+56=", 2) | Out-Null
$d.Content.Find.Execute("23-18=", $true, $true, $false, $false, $false, $true, 1, $false, "88-3=", 2) | Out-Null
$d.Content.Find.Execute("98-93=", $true, $true, $false, $false, $false, $true, 1, $false, "76+2=", 2) | Out-Null
$d.Content.Find.Execute("14+19=", $true, $true, $false, $false, $false, $true, 1, $false, "48+51=", 2) | Out-Null
$d.Content.Find.Execute("11+8=", $true, $true, $false, $false, $false, $true, 1, $false, "63-1=", 2) | Out-Null
$d.Content.Find.Execute("10+3=", $true, $true, $false, $false, $false, $true, 1, $false, "61-23=", 2) | Out-Null
$d.Content.Find.Execute("96-52=", $true, $true, $false, $false, $false, $true, 1, $false, "73-31=", 2) | Out-Null
$d.Content.Find.Execute("20-13=", $true, $true, $false, $false, $false, $true, 1, $false, "29+45=", 2) | Out-Null
$d.Content.Find.Execute("55-0=", $true, $true, $false, $false, $false, $true, 1, $false, "91-42=", 2) | Out-Null
$d.Content.Find.Execute("22-15=", $true, $true, $false, $false, $false, $true, 1, $false, "75-22=", 2) | Out-Null
$d.Content.Find.Execute("51-29=", $true, $true, $false, $false, $false, $true, 1, $false, "13+14=", 2) | Out-Null
$d.Content.Find.Execute("31+30=", $true, $true, $false, $false, $false, $true, 1, $false, "56+20=", 2) | Out-Null
$d.Content.Find.Execute("63-4=", $true, $true, $false, $false, $false, $true, 1, $false, "27+3=", 2) | Out-Null
$d.Content.Find.Execute("23+52=", $true, $true, $false, $false, $false, $true, 1, $false, "72-5=", 2) | Out-Null
$d.Content.Find.Execute("28-7=", $true, $true, $false, $false, $false, $true, 1, $false, "66-40=", 2) | Out-Null
$d.Content.Find.Execute("15-14=", $true, $true, $false, $false, $false, $true, 1, $false, "52-6=", 2) | Out-Null
$d.Content.Find.Execute("42-39=", $true, $true, $false, $false, $false, $true, 1, $false, "89-5=", 2) | Out-Null
$d.Content.Find.Execute("92-85=", $true, $true, $false, $false, $false, $true, 1, $false, "49-16=", 2) | Out-Null
